# SMARTNODES.docx (Romanian) edit.
#
# The "100,000 Smart:" bullet used to read:
#   "... obtained from exchanges such as CryptoBridge, HitBTC."
# with "CryptoBridge" and "HitBTC" rendered via HYPERLINK fields. The
# update drops both hyperlinks (and the "such as " / ", " connective
# text around them) so the sentence simply reads:
#   "... obtained from exchanges."
#
# As a side effect of Word renumbering bookmark ids on any content edit,
# the lone "smarthosting" bookmark's w:id is compacted from 1 to 0.

$nbsp = [char]0x00A0
$d = $word.ActiveDocument

# --- Compute the two plain-text spans to remove, against the pristine
#     (not-yet-edited) document, so the offsets are stable. ---
$rSuchAs = $d.Content
$rSuchAs.Find.Execute(" such as" + $nbsp, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$suchAsStart = $rSuchAs.Start
$suchAsEnd = $rSuchAs.End

$rComma = $d.Content
$rComma.Find.Execute("," + $nbsp, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$commaStart = $rComma.Start
$commaEnd = $rComma.End

# Delete the later span first so the earlier span's offsets stay valid.
$d.Range($commaStart, $commaEnd).Delete()
$d.Range($suchAsStart, $suchAsEnd).Delete()

# --- Remove the two HYPERLINK fields (CryptoBridge, HitBTC). ---
# Field.Delete() removes the begin/instrText/separate/end run set *and*
# the field's display-result run (e.g. "CryptoBridge") together, which is
# exactly what the diff strips.
$removed = 0
$i = 1
while ($i -le $d.Fields.Count -and $removed -lt 2) {
    $f = $d.Fields.Item($i)
    $code = $f.Code.Text
    if ($code -like "*HYPERLINK*crypto-bridge.org*" -or $code -like "*HYPERLINK*hitbtc.com*") {
        $f.Delete()
        $removed = $removed + 1
    } else {
        $i = $i + 1
    }
}
